$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3505
$ws.Range("F5").Value = 3505
$ws.Range("F7").Value = 5037
$ws.Range("F9").Value = 340
$ws.Range("F10").Value = 193
$ws.Range("F11").Value = 675
$ws.Range("F15").Value = 687
$ws.Range("F16").Value = 305
$ws.Range("F19").Value = 154
$ws.Range("F22").Value = 4870
$ws.Range("F26").Value = 5992
$ws.Range("F30").Value = 324
$ws.Range("F36").Value = 971
$ws.Range("F40").Value = 851
$ws.Range("F41").Value = 945

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value = "不可售"
$ws.Range("F4").Value = 22

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1109
$ws.Range("F4").Value = 49

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1109
$ws.Range("F5").Value = 49
$ws.Range("G7").Value = "不可售"
$ws.Range("F8").Value = 3505
$ws.Range("F9").Value = 3505
$ws.Range("F11").Value = 5037
$ws.Range("F13").Value = 340
$ws.Range("F14").Value = 193
$ws.Range("F15").Value = 675
$ws.Range("F18").Value = 687
$ws.Range("F19").Value = 305
$ws.Range("F23").Value = 154
$ws.Range("F26").Value = 4870
$ws.Range("F30").Value = 5992
$ws.Range("F34").Value = 324
$ws.Range("F38").Value = 22
$ws.Range("F41").Value = 971
$ws.Range("F45").Value = 851
$ws.Range("F46").Value = 945
